# Update cached leve-profit figures across all 8 job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Values mirror a scheduled market-data refresh run (Universalis price pulls); see commit message.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4167066.5

$ws.Range("H98").Value = 3234.4167
$ws.Range("I98").Value = 3413.8125
$ws.Range("J98").Value = 1799.25
$ws.Range("K98").Value = 3413.8125
$ws.Range("L98").Value = 1799.25
$ws.Range("M98").Value = -1915.8125
$ws.Range("N98").Value = -4795.25

$ws.Range("H106").Value = 5487.2666
$ws.Range("I106").Value = 5985.852
$ws.Range("K106").Value = 5985.852
$ws.Range("M106").Value = -5354.852

$ws.Range("H116").Value = 3022.647
$ws.Range("I116").Value = 2989
$ws.Range("J116").Value = 3084.3333
$ws.Range("K116").Value = 2989
$ws.Range("L116").Value = 3084.3333
$ws.Range("M116").Value = 453
$ws.Range("N116").Value = -9968.3333

$ws.Range("H121").Value = 888.125
$ws.Range("J121").Value = 890.2273
$ws.Range("L121").Value = 2670.6819
$ws.Range("N121").Value = -6164.6819

$ws.Range("H122").Value = 3234.4167
$ws.Range("I122").Value = 3413.8125
$ws.Range("J122").Value = 1799.25
$ws.Range("K122").Value = 10241.4375
$ws.Range("L122").Value = 5397.75
$ws.Range("M122").Value = -7791.4375
$ws.Range("N122").Value = -10297.75

$ws.Range("H132").Value = 6066329.5
$ws.Range("I132").Value = 7250040.5
$ws.Range("J132").Value = 16252.333
$ws.Range("K132").Value = 21750121.5
$ws.Range("L132").Value = 48756.999
$ws.Range("M132").Value = -21747591.5
$ws.Range("N132").Value = -53816.999

$ws.Range("H135").Value = 3152.1
$ws.Range("I135").Value = 1144.2916
$ws.Range("K135").Value = 10298.6244
$ws.Range("M135").Value = -7763.624400000001

$ws.Range("H137").Value = 1547.5807
$ws.Range("I137").Value = 1527.1875
$ws.Range("K137").Value = 4581.5625
$ws.Range("M137").Value = -2031.5625

$ws.Range("H138").Value = 2992.9358
$ws.Range("I138").Value = 2970.6155
$ws.Range("J138").Value = 2997.4
$ws.Range("K138").Value = 8911.8465
$ws.Range("L138").Value = 8992.200000000001
$ws.Range("M138").Value = -3771.8465
$ws.Range("N138").Value = -19272.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13626.8955
$ws.Range("I32").Value = 11753.609
$ws.Range("J32").Value = 15333.667
$ws.Range("K32").Value = 11753.609
$ws.Range("L32").Value = 15333.667
$ws.Range("M32").Value = -11466.609
$ws.Range("N32").Value = -15907.667

$ws.Range("H45").Value = 1004.5
$ws.Range("I45").Value = 843.8570999999999
$ws.Range("J45").Value = 1379.3334
$ws.Range("K45").Value = 843.8570999999999
$ws.Range("L45").Value = 1379.3334
$ws.Range("M45").Value = -466.8570999999999
$ws.Range("N45").Value = -2133.3334

$ws.Range("H61").Value = 142858370
$ws.Range("I61").Value = 200000800
$ws.Range("J61").Value = 2257
$ws.Range("K61").Value = 200000800
$ws.Range("L61").Value = 2257
$ws.Range("M61").Value = -200000588
$ws.Range("N61").Value = -2681

$ws.Range("H74").Value = 1667.0667
$ws.Range("I74").Value = 896.2727
$ws.Range("J74").Value = 2113.3157
$ws.Range("K74").Value = 896.2727
$ws.Range("L74").Value = 2113.3157
$ws.Range("M74").Value = -22.27269999999999
$ws.Range("N74").Value = -3861.3157

$ws.Range("H77").Value = 1667.0667
$ws.Range("I77").Value = 896.2727
$ws.Range("J77").Value = 2113.3157
$ws.Range("K77").Value = 4481.363499999999
$ws.Range("L77").Value = 10566.5785
$ws.Range("M77").Value = -113.3634999999995
$ws.Range("N77").Value = -19302.5785

$ws.Range("H102").Value = 7578262
$ws.Range("I102").Value = 7578262
$ws.Range("K102").Value = 7578262
$ws.Range("M102").Value = -7576640

$ws.Range("H132").Value = 2334.7908
$ws.Range("I132").Value = 1944.0588
$ws.Range("J132").Value = 3810.889
$ws.Range("K132").Value = 5832.1764
$ws.Range("L132").Value = 11432.667
$ws.Range("M132").Value = -3302.1764
$ws.Range("N132").Value = -16492.667

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H136").Value = 142858370
$ws.Range("I136").Value = 200000800
$ws.Range("J136").Value = 2257
$ws.Range("K136").Value = 600002400
$ws.Range("L136").Value = 6771
$ws.Range("M136").Value = -599999850
$ws.Range("N136").Value = -11871

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 58827540
$ws.Range("I86").Value = 58827540
$ws.Range("K86").Value = 58827540
$ws.Range("M86").Value = -58826417

$ws.Range("H89").Value = 58827540
$ws.Range("I89").Value = 58827540
$ws.Range("K89").Value = 294137700
$ws.Range("M89").Value = -294132084

$ws.Range("H94").Value = 10417652
$ws.Range("J94").Value = 1160
$ws.Range("L94").Value = 1160
$ws.Range("N94").Value = -2062

$ws.Range("H134").Value = 3844.5483
$ws.Range("I134").Value = 639.36664
$ws.Range("K134").Value = 1918.09992
$ws.Range("M134").Value = 616.9000800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1783.7115
$ws.Range("I31").Value = 1760.8431
$ws.Range("K31").Value = 1760.8431
$ws.Range("M31").Value = -1465.8431

$ws.Range("H34").Value = 1783.7115
$ws.Range("I34").Value = 1760.8431
$ws.Range("K34").Value = 1760.8431
$ws.Range("M34").Value = -1558.8431

$ws.Range("H41").Value = 22500
$ws.Range("J41").Value = 22500
$ws.Range("L41").Value = 22500
$ws.Range("N41").Value = -23356

$ws.Range("H58").Value = 5440.12
$ws.Range("I58").Value = 815.8461
$ws.Range("J58").Value = 10449.75
$ws.Range("K58").Value = 815.8461
$ws.Range("L58").Value = 10449.75
$ws.Range("M58").Value = -612.8461
$ws.Range("N58").Value = -10855.75

$ws.Range("H60").Value = 5950
$ws.Range("I60").Value = 5950
$ws.Range("K60").Value = 5950
$ws.Range("M60").Value = -5439

$ws.Range("H74").Value = 32250
$ws.Range("J74").Value = 32250
$ws.Range("L74").Value = 32250
$ws.Range("N74").Value = -33998

$ws.Range("H77").Value = 32250
$ws.Range("J77").Value = 32250
$ws.Range("L77").Value = 96750
$ws.Range("N77").Value = -105486

$ws.Range("H134").Value = 10870951
$ws.Range("I134").Value = 1351.1034
$ws.Range("J134").Value = 29413210
$ws.Range("K134").Value = 4053.3102
$ws.Range("L134").Value = 88239630
$ws.Range("M134").Value = -1518.3102
$ws.Range("N134").Value = -88244700

$ws.Range("H136").Value = 5440.12
$ws.Range("I136").Value = 815.8461
$ws.Range("J136").Value = 10449.75
$ws.Range("K136").Value = 2447.5383
$ws.Range("L136").Value = 31349.25
$ws.Range("M136").Value = 102.4616999999998
$ws.Range("N136").Value = -36449.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 55567770
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 55567770
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 166703310
$ws.Range("M112").ClearContents()
$ws.Range("N112").Value = -166705526

$ws.Range("H113").Value = 702.85297
$ws.Range("I113").Value = 596
$ws.Range("J113").Value = 747.375
$ws.Range("K113").Value = 1788
$ws.Range("L113").Value = 2242.125
$ws.Range("M113").Value = 382
$ws.Range("N113").Value = -6582.125

$ws.Range("H122").Value = 921.2308
$ws.Range("J122").Value = 986.62067
$ws.Range("L122").Value = 8879.58603
$ws.Range("N122").Value = -13779.58603

$ws.Range("H131").Value = 19638226
$ws.Range("J131").Value = 32937.13
$ws.Range("L131").Value = 98811.38999999998
$ws.Range("N131").Value = -108891.39

$ws.Range("H140").Value = 28288.564
$ws.Range("I140").Value = 61210.35
$ws.Range("J140").Value = 2849
$ws.Range("K140").Value = 183631.05
$ws.Range("L140").Value = 8547
$ws.Range("M140").Value = -178451.05
$ws.Range("N140").Value = -18907

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 630
$ws.Range("I97").Value = 630
$ws.Range("K97").Value = 630
$ws.Range("M97").Value = -134

$ws.Range("H132").Value = 10919.667
$ws.Range("I132").Value = 14118.3
$ws.Range("J132").Value = 4522.4
$ws.Range("K132").Value = 42354.89999999999
$ws.Range("L132").Value = 13567.2
$ws.Range("M132").Value = -39824.89999999999
$ws.Range("N132").Value = -18627.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 902.36
$ws.Range("I16").Value = 923.2917
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 923.2917
$ws.Range("L16").Value = 400
$ws.Range("M16").Value = -753.2917
$ws.Range("N16").Value = -740

$ws.Range("H68").Value = 2035
$ws.Range("I68").Value = 2000
$ws.Range("K68").Value = 2000
$ws.Range("M68").Value = -1251

$ws.Range("H71").Value = 2035
$ws.Range("I71").Value = 2000
$ws.Range("K71").Value = 10000
$ws.Range("M71").Value = -6256

$ws.Range("H93").Value = 924.55554
$ws.Range("I93").Value = 952.625
$ws.Range("J93").Value = 700
$ws.Range("K93").Value = 952.625
$ws.Range("L93").Value = 700
$ws.Range("M93").Value = 295.375
$ws.Range("N93").Value = -3196

$ws.Range("H100").Value = 1462.125
$ws.Range("I100").Value = 1139.4
$ws.Range("K100").Value = 1139.4
$ws.Range("M100").Value = -598.4000000000001

$ws.Range("H132").Value = 2634.1765
$ws.Range("I132").Value = 2265.25
$ws.Range("J132").Value = 3519.6
$ws.Range("K132").Value = 6795.75
$ws.Range("L132").Value = 10558.8
$ws.Range("M132").Value = -4265.75
$ws.Range("N132").Value = -15618.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10266.333
$ws.Range("J41").Value = 10266.333
$ws.Range("L41").Value = 10266.333
$ws.Range("N41").Value = -11046.333

$ws.Range("H96").Value = 1410
$ws.Range("J96").Value = 2000
$ws.Range("L96").Value = 2000
$ws.Range("N96").Value = -4746

$ws.Range("H132").Value = 3765.9143
$ws.Range("I132").Value = 3852.2222
$ws.Range("J132").Value = 3474.625
$ws.Range("K132").Value = 11556.6666
$ws.Range("L132").Value = 10423.875
$ws.Range("M132").Value = -9026.6666
$ws.Range("N132").Value = -15483.875

$ws.Range("H136").Value = 952.3125
$ws.Range("I136").Value = 655.9048
$ws.Range("K136").Value = 1967.7144
$ws.Range("M136").Value = 582.2855999999999

